# Revert System_data back to its state before the "oppgave 4" tweak.
#
# - BranchData (sheet1): row 4 (branch 3-4) was a tweak that needs to be
#   removed; deleting it shifts every row below it up by one, which is
#   exactly what the diff shows (rows 5-10 -> rows 4-9).
# - Selection / active-cell bookmarks move back too (G23 -> G14 on
#   BranchData, K19 -> K5 on BusData).
# - BusData (sheet2) column widths for columns J/K go back to explicit
#   (non bestFit) widths.

$wb = $excel.ActiveWorkbook

$wsBranch = $wb.Worksheets.Item("BranchData")
$wsBus    = $wb.Worksheets.Item("BusData")

# --- BranchData: drop the extra row (row 4), rows below shift up ---
$wsBranch.Rows.Item(4).Delete()

# --- BusData: restore the pre-tweak column widths for J and K ---
$wsBus.Columns.Item(10).ColumnWidth = 10.57
$wsBus.Columns.Item(11).ColumnWidth = 13.71

# --- Restore window position bookkeeping ---
$win = $excel.ActiveWindow
$win.Left = -28920
$win.Top = -120

# --- Selections (set BusData first so BranchData ends up the active tab) ---
$wsBus.Activate()
$wsBus.Range("K5").Select()

$wsBranch.Activate()
$wsBranch.Range("G14").Select()
